# edit.ps1 - applies the "edited comments in User manual" change:
#  1. Adds a "Caption" style (Word's standard figure-caption style).
#  2. Removes the old "This is a win horizontally / vertically / diagonally"
#     paragraph (the one using manual tab stops).
#  3. Marks the paragraph containing the three screenshots as keepNext and
#     gives each picture a drop-shadow (outerShdw) effect plus a slightly
#     larger effectExtent to accommodate it; Picture 8 is also resized down
#     a bit.
#  4. Adds a new "Caption"-styled paragraph right after the pictures with
#     the Figure 1/2/3 captions (using SEQ fields for the auto-numbers).

$d = $word.ActiveDocument

# ---- 1. Create the "Caption" style -----------------------------------
$capStyle = $d.Styles.Add("Caption", 1)
$capStyle.BaseStyle = $d.Styles.Item("Normal")
$capStyle.NextParagraphStyle = $d.Styles.Item("Normal")
$capStyle.Priority = 35
$capStyle.UnhideWhenUsed = $true
$capStyle.QuickStyle = $true
$capStyle.NameLocal = "caption"
$capStyle.Font.Italic = $true
$capStyle.Font.ItalicBi = $true
$capStyle.Font.Size = 9
$capStyle.Font.SizeBi = 9
$capStyle.Font.TextColor.ObjectThemeColor = 15
$capStyle.ParagraphFormat.SpaceAfter = 10

# ---- 2/3/4. Replace the "how to win" caption text + the pictures -----
# paragraph with: the pictures paragraph (now keepNext + drop shadows)
# followed by the new Figure 1/2/3 caption paragraph. Locate the two
# paragraphs by their text so the script is resilient to any paragraph
# numbering differences.
$winPara = $null
$picsPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*This is a win horizontally*") {
        $winPara = $para
        $picsPara = $d.Paragraphs.Item($i + 1)
        break
    }
}

if ($winPara -eq $null) {
    Write-Output "ERROR: could not locate the 'This is a win horizontally' paragraph"
} else {
    $target = $d.Range($winPara.Range.Start, $picsPara.Range.End)

    $xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><w:body><w:p w14:paraId="55C01B4F" w14:textId="496743A4" w:rsidR="00871F92" w:rsidRDefault="00FD4885" w:rsidP="00871F92"><w:pPr><w:keepNext/></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="1E615CCD" wp14:editId="43171BD4"><wp:extent cx="1841500" cy="1315731"/><wp:effectExtent l="50800" t="0" r="50800" b="106680"/><wp:docPr id="1" name="Picture 1"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="1" name=""/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId8"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="1892622" cy="1352257"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:effectLst><a:outerShdw blurRad="50800" dist="50800" dir="5400000" algn="ctr" rotWithShape="0"><a:srgbClr val="000000"><a:alpha val="90781"/></a:srgbClr></a:outerShdw></a:effectLst></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="0DF622F0" wp14:editId="15F6F230"><wp:extent cx="1824712" cy="1195754"/><wp:effectExtent l="50800" t="0" r="55245" b="99695"/><wp:docPr id="8" name="Picture 8" descr="A picture containing text&#xA;&#xA;Description automatically generated"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="8" name="Picture 8" descr="A picture containing text&#xA;&#xA;Description automatically generated"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId9"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="1871898" cy="1226676"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:effectLst><a:outerShdw blurRad="50800" dist="50800" dir="5400000" algn="ctr" rotWithShape="0"><a:srgbClr val="000000"/></a:outerShdw></a:effectLst></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r><w:r w:rsidR="00E60C3A"><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="3B87E80C" wp14:editId="3B1B25A9"><wp:extent cx="1793869" cy="1327150"/><wp:effectExtent l="50800" t="0" r="48260" b="95250"/><wp:docPr id="9" name="Picture 9" descr="Text&#xA;&#xA;Description automatically generated with medium confidence"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="9" name="Picture 9" descr="Text&#xA;&#xA;Description automatically generated with medium confidence"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId10"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="1820281" cy="1346690"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:effectLst><a:outerShdw blurRad="50800" dist="50800" dir="5400000" algn="ctr" rotWithShape="0"><a:srgbClr val="000000"/></a:outerShdw></a:effectLst></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Caption"/></w:pPr><w:r><w:t xml:space="preserve">Figure </w:t></w:r><w:fldSimple w:instr=" SEQ Figure \* ARABIC "><w:r><w:rPr><w:noProof/></w:rPr><w:t>1</w:t></w:r></w:fldSimple><w:r><w:t>: This is a win vertically</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve">         </w:t></w:r><w:r><w:t xml:space="preserve">Figure </w:t></w:r><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> SEQ Figure \* ARABIC </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:t>2</w:t></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r><w:r><w:t>: This is a win horizontally</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve">                Figure 3: This is a win diagonally</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

    $target.InsertXML($xml)
    Write-Output "Replacement applied"
}
